# Generate Report for Handoff
# - Status moves from "In Translation" to "Ready for handoff"
# - Timestamps bump forward a few seconds to reflect the new handoff
# - Status/handoff-datetime columns widen to fit the new "Ready for handoff" text

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "Ready for handoff"

# --- Overview sheet ---------------------------------------------------
# E2 = zh-cn status, F2 = de-de status, G2 = Latest HO Xliff Generate Date
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("G2").Value = "2016-09-06 06:45:59"

# --- zh-cn sheet --------------------------------------------------------
# C2 = Status, H2 = Latest Handoff Datetime
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("H2").Value = "2016-09-06 06:45:54"

# --- de-de sheet --------------------------------------------------------
# C2 = Status, H2 = Latest Handoff Datetime
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("H2").Value = "2016-09-06 06:45:59"

# --- Widen the Status / zh-cn / de-de columns to fit the longer text ---
$newWidth = 16.333333333333332

$wsOverview.Columns.Item(5).ColumnWidth = $newWidth  # column E (zh-cn)
$wsOverview.Columns.Item(6).ColumnWidth = $newWidth  # column F (de-de)
$wsZhCn.Columns.Item(3).ColumnWidth = $newWidth       # column C (Status)
$wsDeDe.Columns.Item(3).ColumnWidth = $newWidth       # column C (Status)
